{"js": "// Apply targeted text replacements inside the table cells.\n// Each pair is a unique three-digit x one-digit multiplication problem\n// of the form 'NNN\u00d7N=NNNN' that appears exactly once in the document body.\nconst replacements = [\n  [\"617\u00d72=1234\", \"769\u00d72=1538\"],\n  [\"646\u00d77=4522\", \"154\u00d75=770\"],\n  [\"456\u00d79=4104\", \"208\u00d75=1040\"],\n  [\"382\u00d78=3056\", \"422\u00d78=3376\"],\n  [\"367\u00d76=2202\", \"765\u00d77=5355\"],\n  [\"463\u00d75=2315\", \"222\u00d74=888\"],\n  [\"868\u00d79=7812\", \"559\u00d79=5031\"],\n  [\"842\u00d72=1684\", \"269\u00d76=1614\"],\n  [\"985\u00d76=5910\", \"807\u00d78=6456\"],\n  [\"298\u00d78=2384\", \"400\u00d75=2000\"],\n  [\"164\u00d78=1312\", \"413\u00d72=826\"],\n  [\"395\u00d73=1185\", \"407\u00d78=3256\"],\n  [\"656\u00d74=2624\", \"251\u00d72=502\"],\n  [\"601\u00d75=3005\", \"440\u00d73=1320\"],\n  [\"265\u00d76=1590\", \"727\u00d74=2908\"],\n  [\"972\u00d79=8748\", \"321\u00d76=1926\"],\n  [\"261\u00d72=522\", \"224\u00d77=1568\"],\n  [\"503\u00d74=2012\", \"393\u00d74=1572\"],\n  [\"586\u00d79=5274\", \"678\u00d72=1356\"],\n  [\"995\u00d77=6965\", \"181\u00d74=724\"],\n  [\"548\u00d78=4384\", \"349\u00d78=2792\"],\n  [\"531\u00d73=1593\", \"125\u00d79=1125\"],\n  [\"961\u00d79=8649\", \"194\u00d73=582\"],\n  [\"574\u00d77=4018\", \"320\u00d77=2240\"],\n  [\"302\u00d75=1510\", \"980\u00d75=4900\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply targeted text replacements inside the table cells.\n# Each pair is a unique three-digit x one-digit multiplication problem\n# of the form 'NNN x N=NNNN' that appears exactly once in the document.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"617\u00d72=1234\", \"769\u00d72=1538\")\n    ,@(\"646\u00d77=4522\", \"154\u00d75=770\")\n    ,@(\"456\u00d79=4104\", \"208\u00d75=1040\")\n    ,@(\"382\u00d78=3056\", \"422\u00d78=3376\")\n    ,@(\"367\u00d76=2202\", \"765\u00d77=5355\")\n    ,@(\"463\u00d75=2315\", \"222\u00d74=888\")\n    ,@(\"868\u00d79=7812\", \"559\u00d79=5031\")\n    ,@(\"842\u00d72=1684\", \"269\u00d76=1614\")\n    ,@(\"985\u00d76=5910\", \"807\u00d78=6456\")\n    ,@(\"298\u00d78=2384\", \"400\u00d75=2000\")\n    ,@(\"164\u00d78=1312\", \"413\u00d72=826\")\n    ,@(\"395\u00d73=1185\", \"407\u00d78=3256\")\n    ,@(\"656\u00d74=2624\", \"251\u00d72=502\")\n    ,@(\"601\u00d75=3005\", \"440\u00d73=1320\")\n    ,@(\"265\u00d76=1590\", \"727\u00d74=2908\")\n    ,@(\"972\u00d79=8748\", \"321\u00d76=1926\")\n    ,@(\"261\u00d72=522\", \"224\u00d77=1568\")\n    ,@(\"503\u00d74=2012\", \"393\u00d74=1572\")\n    ,@(\"586\u00d79=5274\", \"678\u00d72=1356\")\n    ,@(\"995\u00d77=6965\", \"181\u00d74=724\")\n    ,@(\"548\u00d78=4384\", \"349\u00d78=2792\")\n    ,@(\"531\u00d73=1593\", \"125\u00d79=1125\")\n    ,@(\"961\u00d79=8649\", \"194\u00d73=582\")\n    ,@(\"574\u00d77=4018\", \"320\u00d77=2240\")\n    ,@(\"302\u00d75=1510\", \"980\u00d75=4900\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute(\n        $oldText,   # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
